$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row values (F1, G1)
$ws.Range("F1").Value = "pidan add prefetch in pack"
$ws.Range("G1").Value = "pidan tune mc nc kc"

# New data row values (F2, G2)
$ws.Range("F2").Value = "62(81.5%)"
$ws.Range("G2").Value = "63.31(83%)"

# Column widths: F matches E's width (~26.16), G gets width ~19
$ws.Columns.Item(6).ColumnWidth = 25.45
$ws.Columns.Item(7).ColumnWidth = 18.25

# Update selection to D5 (matches final cursor position in the diff)
[void]$ws.Range("D5").Select()
